# Auto-generated Excel COM-interop script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row reorder / content swap for rows 13-16 (Litecoin/Polkadot/TRON/Chainlink) ---
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"

# --- Price column (D): force text storage so values like "26.943.53" / "93.00" keep their
#     original string formatting instead of being auto-converted to numbers, then restore the
#     default "Normal" style so no stray style index is left on the cells. ---
$dCells = @("D2","D3","D5","D7","D8","D9","D10","D12","D13","D14","D15","D16","D18","D20","D21","D22","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "26.943.53"
$ws.Range("D3").Value = "1.819.75"
$ws.Range("D5").Value = "309.93"
$ws.Range("D7").Value = "0.4686"
$ws.Range("D8").Value = "0.3703"
$ws.Range("D9").Value = "0.07394"
$ws.Range("D10").Value = "0.8747"
$ws.Range("D12").Value = "1.849.41"
$ws.Range("D13").Value = "5.376"
$ws.Range("D14").Value = "93.00"
$ws.Range("D15").Value = "6.525"
$ws.Range("D16").Value = "0.07070"
$ws.Range("D18").Value = "0.000008733"
$ws.Range("D20").Value = "14.80"
$ws.Range("D21").Value = "26.964.00"
$ws.Range("D22").Value = "5.332"
$ws.Range("D24").Value = "2.046.10"
$ws.Range("D25").Value = "1.902"
$ws.Range("D26").Value = "151.74"
$ws.Range("D27").Value = "2.216"
$ws.Range("D28").Value = "18.46"
$ws.Range("D29").Value = "5.335"
$ws.Range("D30").Value = "115.79"
$ws.Range("D31").Value = "0.08942"
$ws.Range("D32").Value = "0.7700"
$ws.Range("D33").Value = "1.171"
$ws.Range("D34").Value = "4.503"
$ws.Range("D35").Value = "2.910"
$ws.Range("D36").Value = "0.9999"
$ws.Range("D37").Value = "1.086"
$ws.Range("D38").Value = "0.01966"
$ws.Range("D39").Value = "0.05294"
$ws.Range("D40").Value = "7.324"
$ws.Range("D42").Value = "0.5364"
$ws.Range("D43").Value = "2.380"
$ws.Range("D44").Value = "0.1673"
$ws.Range("D45").Value = "8.482"
$ws.Range("D46").Value = "0.4968"
$ws.Range("D47").Value = "10.44"
$ws.Range("D48").Value = "1.679"
$ws.Range("D49").Value = "0.9997"
$ws.Range("D50").Value = "103.52"
$ws.Range("D51").Value = "0.06299"

foreach ($addr in $dCells) { $ws.Range($addr).Style = "Normal" }

# --- Volume(1h) column (E): plain percentage strings, safe to assign directly ---
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("E7").Value = "  +2.45%  "
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("E10").Value = "  +2.15%  "
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("E12").Value = "  +2.81%  "
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("E14").Value = "  +2.14%  "
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("E22").Value = "  +0.94%  "
$ws.Range("E23").Value = "  -1.78%  "
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("E26").Value = "  +0.34%  "
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("E29").Value = "  +2.09%  "
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("E35").Value = "  +0.39%  "
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("E37").Value = "  -3.19%  "
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("E39").Value = "  +1.62%  "
$ws.Range("E40").Value = "  +1.88%  "
$ws.Range("E41").Value = "  +2.04%  "
$ws.Range("E42").Value = "  +2.33%  "
$ws.Range("E43").Value = "  +2.89%  "
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("E51").Value = "  -0.13%  "
